# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted for Perejil (Vega Central
# Mapocho de Santiago) at sheet row 142, pushing every subsequent record
# down by one row (142-214 -> 143-215) and extending the used range to
# A1:R215.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 142; Excel shifts rows 142:214 down
# to 143:215 and the sheet dimension grows to R215 automatically.
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new record's data.
$ws.Cells.Item(142, 1).Value2 = 9
$ws.Cells.Item(142, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(142, 3).Value2 = "Metropolitana"
$ws.Cells.Item(142, 4).Value2 = 44466
$ws.Cells.Item(142, 5).Value2 = 13
$ws.Cells.Item(142, 6).Value2 = 100112044
$ws.Cells.Item(142, 7).Value2 = "Perejil"
$ws.Cells.Item(142, 8).Value2 = "Sin especificar"
$ws.Cells.Item(142, 9).Value2 = "Primera"
$ws.Cells.Item(142, 10).Value2 = 106
$ws.Cells.Item(142, 11).Value2 = 12000
$ws.Cells.Item(142, 12).Value2 = 14000
$ws.Cells.Item(142, 13).Value2 = 13000
$ws.Cells.Item(142, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(142, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(142, 16).Value2 = 4333
$ws.Cells.Item(142, 17).Value2 = 3
$ws.Cells.Item(142, 18).Value2 = "Hortaliza"

# Match the date formatting used by the rest of column D.
$ws.Cells.Item(142, 4).NumberFormat = $ws.Cells.Item(143, 4).NumberFormat
